$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number must be
# forced to Text format first, otherwise Excel auto-converts the string
# to a numeric value (losing the original text formatting / trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "26.926.82"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.554.61"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "207.21"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "21.72"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").Value = "0.0587"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.775.55"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "1.555.28"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "3.72"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "61.81"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "26.907.71"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "215.37"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "7.23"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "152.30"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "14.88"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "0.0463"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "1.416.97"
$ws.Range("E33").Value = "  +4.55%  "
$ws.Range("E34").Value = "  +3.49%  "
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +4.42%  "
$ws.Range("D36").Value = "0.961"
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "2.28"
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("D45").Value = "63.63"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "1.690.21"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("D50").Value = "0.0959"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +0.24%  "
